# Re-saving-in-newer-Excel style update:
#  - default/base font switches from the Chinese "SimSun" (宋体) face to Calibri
#  - worksheet is renamed from the generic "Sheet1" to "Properties"
#  - the active selection on the sheet moves to I22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the workbook's base/Normal style font (used by every cell that has
# no explicit formatting) from SimSun to Calibri.
$normalStyle = $wb.Styles.Item(1)
$normalStyle.Font.Name = "Calibri"
$normalStyle.Name = "Normal"

# Rename the worksheet.
$ws.Name = "Properties"

# Move the selection/active cell.
[void]$ws.Range("I22").Select()
